# Rename the two earliest sprint labels to include their date ranges.
# (Excel will append the new/changed shared-string text at the end of the
# shared-strings table and re-point the B6/B7 cells at it; "Sprint 3" and
# "Sprint 4" -- still referenced unchanged by B8/B9 -- simply shift to fill
# the now-unused slots, exactly matching the committed workbook.)
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Planilha1")

$ws.Range("B6").Value = "Sprint 1 10/03 - 17/03"
$ws.Range("B7").Value = "Sprint 2 17/03 - 24/03"

# Column B needs to widen so the longer sprint labels are readable.
$ws.Columns.Item(2).ColumnWidth = 25.65

# Leave the selection where the edit happened, like a user would.
[void]$ws.Range("B7").Select()
